$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.9
$ws.Range("I2").Value = 2.63
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("U2").Value = 2.2
$ws.Range("V2").Value = 1.62
$ws.Range("X2").Value = 13
$ws.Range("AO2").Value = 19
$ws.Range("AQ2").Value = 67
$ws.Range("AS2").Value = 351
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("G4").Value = 4.5
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 1.95
$ws.Range("J4").Value = 5
$ws.Range("L4").Value = 2.75
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("W4").Value = 9.5
$ws.Range("X4").Value = 21
$ws.Range("Y4").Value = 17
$ws.Range("Z4").Value = 51
$ws.Range("AN4").Value = 6
$ws.Range("AO4").Value = 26
$ws.Range("AQ4").Value = 101
$ws.Range("I5").Value = 13
$ws.Range("K5").Value = 2.3
$ws.Range("L5").Value = 11
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("O5").Value = 1.3
$ws.Range("P5").Value = 3.4
$ws.Range("Q5").Value = 2.03
$ws.Range("R5").Value = 1.83
$ws.Range("S5").Value = 1.4
$ws.Range("T5").Value = 2.75
$ws.Range("U5").Value = 2.63
$ws.Range("V5").Value = 1.44
$ws.Range("X5").Value = 5
$ws.Range("Y5").Value = 10
$ws.Range("Z5").Value = 7
$ws.Range("AC5").Value = 9
$ws.Range("AE5").Value = 34
$ws.Range("AF5").Value = 126
$ws.Range("AL5").Value = 101
$ws.Range("AM5").Value = 101
$ws.Range("AN5").Value = 3
$ws.Range("AP5").Value = 23
$ws.Range("AR5").Value = 51
$ws.Range("AS5").Value = 251
$ws.Range("AT5").Value = 2.75
$ws.Range("AU5").Value = 12
$ws.Range("AV5").Value = 101
$ws.Range("AW5").Value = 11
$ws.Range("AX5").Value = 51
$ws.Range("AY5").Value = 51
$ws.Range("AZ5").Value = 351
$ws.Range("G6").Value = 1.5
$ws.Range("H6").Value = 4
$ws.Range("J6").Value = 2.1
$ws.Range("U6").Value = 2.1
$ws.Range("V6").Value = 1.67
$ws.Range("X6").Value = 6.5
$ws.Range("Y6").Value = 9
$ws.Range("Z6").Value = 10
$ws.Range("AH6").Value = 15
$ws.Range("AI6").Value = 34
$ws.Range("AJ6").Value = 21
$ws.Range("AS6").Value = 151
$ws.Range("AU6").Value = 9.5
$ws.Range("O7").Value = 1.36
$ws.Range("P7").Value = 3
$ws.Range("Q7").Value = 2.1
$ws.Range("R7").Value = 1.7
$ws.Range("K8").Value = 2.4
$ws.Range("O8").Value = 1.22
$ws.Range("P8").Value = 4
$ws.Range("Q8").Value = 1.8
$ws.Range("R8").Value = 2
$ws.Range("U8").Value = 2
$ws.Range("V8").Value = 1.73
$ws.Range("AK8").Value = 101
$ws.Range("G9").Value = 1.4
$ws.Range("H9").Value = 4.33
$ws.Range("I9").Value = 8.5
$ws.Range("J9").Value = 1.91
$ws.Range("L9").Value = 7.5
$ws.Range("M9").Value = 1.05
$ws.Range("N9").Value = 11
$ws.Range("U9").Value = 2.1
$ws.Range("V9").Value = 1.67
$ws.Range("X9").Value = 6.5
$ws.Range("Z9").Value = 9
$ws.Range("AB9").Value = 29
$ws.Range("AC9").Value = 11
$ws.Range("AD9").Value = 8.5
$ws.Range("AE9").Value = 21
$ws.Range("AF9").Value = 67
$ws.Range("AJ9").Value = 23
$ws.Range("AL9").Value = 51
$ws.Range("AM9").Value = 51
$ws.Range("AN9").Value = 3.25
$ws.Range("AU9").Value = 9.5
$ws.Range("AW9").Value = 8.5
$ws.Range("G10").Value = 1.7
$ws.Range("I10").Value = 5.5
$ws.Range("J10").Value = 2.38
$ws.Range("L10").Value = 5.5
$ws.Range("M10").Value = 1.06
$ws.Range("N10").Value = 10
$ws.Range("O10").Value = 1.33
$ws.Range("P10").Value = 3.25
$ws.Range("Q10").Value = 2.05
$ws.Range("R10").Value = 1.75
$ws.Range("S10").Value = 1.44
$ws.Range("T10").Value = 2.63
$ws.Range("U10").Value = 1.91
$ws.Range("V10").Value = 1.8
$ws.Range("W10").Value = 6.5
$ws.Range("Y10").Value = 8.5
$ws.Range("AB10").Value = 29
$ws.Range("AC10").Value = 8.5
$ws.Range("AE10").Value = 17
$ws.Range("AF10").Value = 51
$ws.Range("AG10").Value = 351
$ws.Range("AH10").Value = 13
$ws.Range("AN10").Value = 3.6
$ws.Range("AO10").Value = 9
$ws.Range("AQ10").Value = 29
$ws.Range("AT10").Value = 2.63
$ws.Range("AU10").Value = 9
$ws.Range("AV10").Value = 67
$ws.Range("AW10").Value = 6.5
$ws.Range("AX10").Value = 29
$ws.Range("AZ10").Value = 101
$ws.Range("BA10").Value = 126
$ws.Range("BB10").Value = 301
$ws.Range("Q11").Value = 2.08
$ws.Range("R11").Value = 1.73
$ws.Range("G12").Value = 2.77
$ws.Range("H12").Value = 2.6
$ws.Range("I12").Value = 2.9
$ws.Range("J12").Value = 3.4
$ws.Range("K12").Value = 1.82
$ws.Range("L12").Value = 3.6
$ws.Range("S12").Value = 1.53
$ws.Range("T12").Value = 2.18
$ws.Range("W12").Value = 6.6
$ws.Range("X12").Value = 13
$ws.Range("Y12").Value = 10.5
$ws.Range("AA12").Value = 29
$ws.Range("AC12").Value = 5.8
$ws.Range("AD12").Value = 5.2
$ws.Range("AE12").Value = 16
$ws.Range("AF12").Value = 100
$ws.Range("AH12").Value = 6.7
$ws.Range("AI12").Value = 13.5
$ws.Range("AJ12").Value = 10.75
$ws.Range("AK12").Value = 40
$ws.Range("AL12").Value = 32
$ws.Range("AN12").Value = 4.45
$ws.Range("AO12").Value = 16
$ws.Range("AQ12").Value = 80
$ws.Range("AU12").Value = 7
$ws.Range("AV12").Value = 75
$ws.Range("AW12").Value = 4.6
$ws.Range("AX12").Value = 17
$ws.Range("AY12").Value = 26
$ws.Range("AZ12").Value = 90
$ws.Range("AZ9").Value = 151
$ws.Range("BA9").Value = 201
